# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the Slide Master ("Integral" palette)
#   ppt/theme/theme2.xml  -> bound to the Notes Master  ("Office Theme" palette)
#
# The target edit swaps the two palettes: the Slide Master's theme should
# become the stock "Office" color scheme (what the Notes Master used to
# have), and the Notes Master's theme should become the "Integral" colors
# (what the Slide Master used to have).
#
# Apply the reachable half of that swap through the PowerPoint object model:
# push the stock Office theme color values onto the active/Slide Master
# theme's ThemeColorScheme (the twelve dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink slots), in the same order PowerPoint exposes them.

function Convert-HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Stock "Office" theme color scheme, in ThemeColorScheme.Item(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink.
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slideMasterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $slideMasterScheme.Count; $i++) {
    $slideMasterScheme.Item($i).RGB = Convert-HexToComRGB $officeThemeHex[$i - 1]
}

# Mirror the same call against the Notes Master's theme so that, on hosts
# where the Notes Master owns an independently addressable theme part
# (ppt/theme/theme2.xml), it picks up the palette the Slide Master used to
# have -- completing the swap described by the diff.
$notesMasterScheme = $p.NotesMaster.Theme.ThemeColorScheme
$integralThemeHex = @(
    "000000",
    "FFFFFF",
    "455F51",
    "E3DED1",
    "99CB38",
    "63A537",
    "E6D024",
    "CC9700",
    "4EB3CF",
    "378DA6",
    "6B9F25",
    "B26B02"
)
for ($i = 1; $i -le $notesMasterScheme.Count; $i++) {
    $notesMasterScheme.Item($i).RGB = Convert-HexToComRGB $integralThemeHex[$i - 1]
}

# Re-assert the Slide Master's Office palette last so that, on hosts (like
# this sandbox) where Theme/ThemeColorScheme always resolves back to the
# single primary theme part regardless of which master it was reached
# through, the Slide Master's intended "Office Theme" colors are the ones
# that stick.
for ($i = 1; $i -le $slideMasterScheme.Count; $i++) {
    $slideMasterScheme.Item($i).RGB = Convert-HexToComRGB $officeThemeHex[$i - 1]
}
